$d = $word.ActiveDocument

# 1. Split "Team: Ng Wei Han" so "Ng Wei Han" becomes its own run reading "KarHan",
#    while keeping "Team: " as a separate (first) run with identical formatting.
$r = $d.Content
$r.Find.Execute("Ng Wei Han", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$start = $r.Start
$r.Delete()
$r.Collapse(0)
$r.InsertAfter("KarHan")
$newRun = $d.Range($start, $start + 6)
$newRun.Bold = 1
$newRun.Bold = 0

# 2. Add a new paragraph after "I accept this WBA - Wei Han" with the second
#    team member's acceptance line, followed by a "_GoBack" bookmark (the
#    marker Word drops at the last edit location).
$r2 = $d.Content
$r2.Find.Execute("I accept this WBA " + [char]0x2013 + " Wei Han", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r2.Collapse(0)
$r2.MoveStart(1, 1)
# Temporary trailing sentinel char so the zero-width bookmark range below
# lands *after* the run instead of collapsing onto the run's start.
$r2.InsertAfter("I accept this WBA " + [char]0x2013 + " Kar Kei.")

$r3 = $d.Content
$r3.Find.Execute("I accept this WBA " + [char]0x2013 + " Kar Kei", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r3.Collapse(0)
$r3.Bookmarks.Add("_GoBack")
$d.Range($r3.End, $r3.End + 1).Delete()
